$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.68%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "28.37"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-3.16%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.291"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.93%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.46%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.208"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.84%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8500"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.00%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8837"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.52%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1396"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.18%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07099"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.40%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03146"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "4.57%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09226"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.69%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001542"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.11%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0005978"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-94.18%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005936"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.00%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.496"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.20%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-3.69%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03329"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.25%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.06%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.513"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.89%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.82%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.08%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001225"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.02%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004154"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-16.97%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001199"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.84%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001445"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.72%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1066"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.53%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-35.05%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002199"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-9.43%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009479"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "0.27%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005275"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.39%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.01%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-0.34%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.01%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.01%"
